# Fill in the weekly progress report fields.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6
$ws.Range("B3").Value = Get-Date -Year 2025 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("B4").Value = Get-Date -Year 2025 -Month 11 -Day 16 -Hour 0 -Minute 0 -Second 0
$ws.Range("B5").Value = "Lại Việt Anh"
$ws.Range("B6").Value = 179066
$ws.Range("B7").Value = "Xây dựng mô hình giám sát và điều khiển nhà thông minh"
$ws.Range("B8").Value = "thiết kế hệ thống"
$ws.Range("B9").Value = "thiết kế hệ thống"

$ws.Range("K6").Select()
